# Insert a new parameter row for "chemical_recycling_pyrolysis" right after the
# existing "chemical_recycling_gasification" row (row 9), pushing every row
# below it down by one. Matches the commit: "added pyrolysis and additional figures".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 currently holds "fossil_routes"; insert a fresh row above it so that
# row becomes row 11 and the new row takes its place as row 10.
$ws.Rows.Item(10).Insert()

$ws.Range("A10").Value = "chemical_recycling_pyrolysis"
$ws.Range("B10").Value = $true
